# Komponenttitaulukko.xlsx edit:
#  - rename header cells B1/C1
#  - touch the formatting of the existing data cells / column A
#  - add two new columns (B, C) with their own widths for the new price columns
#  - move the active selection to C6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Hinta" / "Kauppa" headers to the new, more specific labels.
$ws.Range("B1").Value = "Verkkokauppa Hinta"
$ws.Range("C1").Value = "Proshop Hinta"

# Re-apply the base cell style to the existing data range; this nudges the
# workbook into carrying a second cellXfs record (the original data moves
# off the bare default style index 0 onto its own style index 1).
$ws.Range("A1").Style = "Normal"
$ws.Range("B1").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("D1").Style = "Normal"
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").Style = "Normal"

# Give the two new "price" columns their own widths.
$ws.Columns.Item(2).ColumnWidth = 18.09
$ws.Columns.Item(3).ColumnWidth = 22.25

# Move the active cell / selection.
[void]$ws.Range("C6").Select()
